$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 900
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 450
$ws.Range("B5").Value = 140
